$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 9
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 7.8
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 3.6
$ws.Range("Q2").Value = 1.64
$ws.Range("R2").Value = 1.43
$ws.Range("S2").Value = 2.6
$ws.Range("T2").Value = 2.34
$ws.Range("U2").Value = 1.52
$ws.Range("Y2").Value = 60
$ws.Range("AC2").Value = 19.5
$ws.Range("AD2").Value = 95
$ws.Range("AJ2").Value = 10
$ws.Range("AK2").Value = 19.5
$ws.Range("H3").Value = 2.46
$ws.Range("K3").Value = 3.6
$ws.Range("P3").Value = 1.78
$ws.Range("T3").Value = 1.81
$ws.Range("F4").Value = 15.5
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 1.24
$ws.Range("I4").Value = 1.31
$ws.Range("J4").Value = 6
$ws.Range("K4").Value = 7.2
$ws.Range("N4").Value = 2.34
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1.58
$ws.Range("R4").Value = 1.34
$ws.Range("S4").Value = 2.24
$ws.Range("T4").Value = 2.28
$ws.Range("V4").Value = 4.3
$ws.Range("W4").Value = 1.05
$ws.Range("X4").Value = 30
$ws.Range("Z4").Value = 9.800000000000001
$ws.Range("AA4").Value = 12.5
$ws.Range("AD4").Value = 16
$ws.Range("AE4").Value = 22
$ws.Range("AO4").Value = 6.8
$ws.Range("F5").Value = 7.2
$ws.Range("H5").Value = 1.02
$ws.Range("I5").Value = 1.54
$ws.Range("J5").Value = 4.7
$ws.Range("K5").Value = 950
$ws.Range("N5").Value = 2.38
$ws.Range("O5").Value = 1.17
$ws.Range("P5").Value = 2.38
$ws.Range("Q5").Value = 1.17
$ws.Range("R5").Value = 1.54
$ws.Range("S5").Value = 1.17
$ws.Range("V5").Value = 2.84
$ws.Range("G6").Value = 2.8
$ws.Range("K6").Value = 3.95
$ws.Range("Q6").Value = 1.86
$ws.Range("W6").Value = 1.55
$ws.Range("I7").Value = 3.3
$ws.Range("K7").Value = 3.8
$ws.Range("Q7").Value = 1.9
$ws.Range("T7").Value = 1.73
$ws.Range("U7").Value = 2.14
$ws.Range("V7").Value = 1.43
$ws.Range("X7").Value = 17.5
$ws.Range("AN7").Value = 28
$ws.Range("F8").Value = 8.4
$ws.Range("G8").Value = 8.800000000000001
$ws.Range("H8").Value = 1.42
$ws.Range("J8").Value = 5.5
$ws.Range("K8").Value = 5.6
$ws.Range("P8").Value = 2.54
$ws.Range("R8").Value = 1.6
$ws.Range("T8").Value = 1.9
$ws.Range("U8").Value = 2.06
$ws.Range("X8").Value = 25
$ws.Range("AA8").Value = 11.5
$ws.Range("AD8").Value = 10
$ws.Range("AL8").Value = 95
$ws.Range("AN8").Value = 130
$ws.Range("F9").Value = 4.7
$ws.Range("H9").Value = 1.74
$ws.Range("J9").Value = 4.1
$ws.Range("K9").Value = 4.4
$ws.Range("L9").Value = 1.28
$ws.Range("H10").Value = 1.87
$ws.Range("L10").Value = 1.35
$ws.Range("O10").Value = 1.26
$ws.Range("Z10").Value = 12
$ws.Range("F12").Value = 2.36
$ws.Range("I12").Value = 3.5
$ws.Range("N12").Value = 3.7
$ws.Range("AE12").Value = 40
$ws.Range("M13").Value = 1.04
$ws.Range("S13").Value = 2.92
$ws.Range("X13").Value = 22
$ws.Range("AE13").Value = 14
$ws.Range("AH13").Value = 34
$ws.Range("W14").Value = 1.57
$ws.Range("I15").Value = 1.94
$ws.Range("Q15").Value = 1.69
$ws.Range("V15").Value = 2.06
$ws.Range("H16").Value = 1.78
$ws.Range("I16").Value = 1.79
$ws.Range("P16").Value = 2.64
$ws.Range("U16").Value = 2.58
$ws.Range("V16").Value = 2.26
$ws.Range("X16").Value = 23
$ws.Range("AA16").Value = 20
$ws.Range("F17").Value = 2.04
$ws.Range("H17").Value = 4.1
$ws.Range("I17").Value = 4.6
$ws.Range("J17").Value = 3.3
$ws.Range("L17").Value = 1.48
$ws.Range("N17").Value = 3.05
$ws.Range("P17").Value = 1.62
